$d = $word.ActiveDocument

# Trim the trailing clause from the SmartCash mining paragraph:
# "...created for quite some time, until Smartcash reaches a considerable market cap."
# becomes
# "...created for quite some time."
$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()
$find.Execute(
    ", until Smartcash reaches a considerable market cap.",  # FindText
    $true,    # MatchCase
    $false,   # MatchWholeWord
    $false,   # MatchWildcards
    $false,   # MatchSoundsLike
    $false,   # MatchAllWordForms
    $true,    # Forward
    1,        # Wrap (wdFindContinue)
    $false,   # Format
    ".",      # ReplaceWith
    2         # Replace (wdReplaceAll)
)

# Recreate the "exchanges" bookmark in place (no textual effect, kept for
# completeness / parity with the underlying re-serialization of the bookmark
# that happens when the source document is resaved).
if ($d.Bookmarks.Exists("exchanges")) {
    $bm = $d.Bookmarks("exchanges")
    $bmRange = $bm.Range
    $bm.Delete()
    $d.Bookmarks.Add("exchanges", $bmRange)
}

$d.Save()
